# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF holds the game date, stored as text (e.g. "5-19-2012-13").
# It needs to become the correct ISO-style text date "2013-05-19".
#
# NOTE: assigning an ISO-looking string straight to .Value lets Excel's
# normal "looks like a date" autodetection kick in and silently turn the
# cell into a date serial number. Pre-formatting the range as Text keeps
# the corrected value stored as a literal string, same as the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "5-19-2012-13"
$newDate = "2013-05-19"

$dataRange = $ws.Range("BF2:BF31")
$dataRange.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)  # column BF
    if ($cell.Value() -eq $oldDate) {
        $cell.Value = $newDate
    }
}
